# edit.ps1 - applies the changes described by the target diff to the
# "Aufgabe_Excel_04.docx" document via the Word COM-interop object model.

$d = $word.ActiveDocument

function New-XmlPackage($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($para, $pPrXml, $runsXml) {
    $body = '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body>'
    $xml = New-XmlPackage $body
    $para.Range.InsertXML($xml)
}

function Get-ParaByText($doc, $matchText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $matchText) {
            return $p
        }
    }
    return $null
}

function Get-ParaIndexByText($doc, $matchText) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like $matchText) {
            return $idx
        }
    }
    return -1
}

function Insert-ParagraphBefore($doc, $beforeIndex, $pPrXml, $runsXml) {
    $target = $doc.Paragraphs.Item($beforeIndex)
    $target.Range.InsertParagraphBefore()
    $newPara = $doc.Paragraphs.Item($beforeIndex)
    Set-ParagraphXml $newPara $pPrXml $runsXml
}

# ---------------------------------------------------------------------
# Change 1: "Die Formel, SVERWEIS(), ..." -> split out ZÄHLENWENNS() with
# proofErr markers and drop SVERWEIS()/ZÄHLENWENN()/SUMMEWENN().
# ---------------------------------------------------------------------
$pPr12 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr></w:pPr>'
$p = Get-ParaByText $d "*Die Formel, SVERWEIS*"
$runs = '<w:r><w:t xml:space="preserve">Die Formel, </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>ZÄHLENWENNS(</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>), SUMMEWENNS()</w:t></w:r>' + `
        '<w:r><w:t>, SUMMENPRODUKT()</w:t></w:r>'
Set-ParagraphXml $p $pPr12 $runs

# ---------------------------------------------------------------------
# Change 2: "Bes-000" run split into "Bes-00" + "1"
# ---------------------------------------------------------------------
$pPr17_0 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr>'
$p = Get-ParaByText $d "*In der Spalte A stehen*"
$runs = '<w:r><w:t xml:space="preserve">In der Spalte A stehen </w:t></w:r>' + `
        '<w:r><w:t>Bestellnummern</w:t></w:r>' + `
        '<w:r><w:t>. Passen Sie das Zahlenformat so an, dass statt &quot;1&quot; in Zelle A2 &quot;</w:t></w:r>' + `
        '<w:r><w:t>Bes</w:t></w:r>' + `
        '<w:r><w:t>-001&quot; angezeigt werden. (Wichtig: In der Zelle steht eine &quot;1&quot;, nur durch das Zahlenformat wird daraus eben &quot;</w:t></w:r>' + `
        '<w:r><w:t>Bes-00</w:t></w:r>' + `
        '<w:r><w:t>1</w:t></w:r>' + `
        '<w:r><w:t>&quot;. Wenn die Zelle angeklickt ist, dann sehen Sie in der Bearbeitunsleiste nur die &quot;1&quot;.</w:t></w:r>'
Set-ParagraphXml $p $pPr17_0 $runs

# ---------------------------------------------------------------------
# Change 3: "In den Zellen B2:O2 bzw. ..." -> split cell-range list
# ---------------------------------------------------------------------
$p = Get-ParaByText $d "*In den Zellen B2:O2*"
$runs = '<w:r><w:t xml:space="preserve">In den Zellen </w:t></w:r>' + `
        '<w:r><w:t>A</w:t></w:r>' + `
        '<w:r><w:t>2:</w:t></w:r>' + `
        '<w:r><w:t>E</w:t></w:r>' + `
        '<w:r><w:t>2</w:t></w:r>' + `
        '<w:r><w:t>, G2:H2, J2:M2 und</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> O2:O4 verwenden Sie bitte das folgende Format:</w:t></w:r>'
Set-ParagraphXml $p $pPr17_0 $runs

# ---------------------------------------------------------------------
# Changes 4-6: Spalte K/L/M block gets a new "Spalte J" item inserted at
# the top, shifting the existing content down by one entry.
# ---------------------------------------------------------------------
$pPr17_1 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="17"/></w:numPr></w:pPr>'

$ellipsis = [char]8230

# "Spalte K: ..." paragraph becomes "Spalte J: Alle Produkte (Prod-01 - …)"
$p = Get-ParaByText $d "*Spalte K: In wievielen Bestellungen das Produkt bestellt wird*"
$runs = '<w:r><w:t>Spalte J: Alle Produkte (Prod-01 - ' + $ellipsis + ')</w:t></w:r>'
Set-ParagraphXml $p $pPr17_1 $runs

# "Spalte L: ..." paragraph (has lastRenderedPageBreak) becomes "Spalte K: ..."
$p = Get-ParaByText $d "*Spalte L: Wieviele Einheiten vom Produkt insgesamt bestellt wurden*"
$runs = '<w:r><w:lastRenderedPageBreak/><w:t>Spalte K: In wievielen Bestellungen das Produkt bestellt wird</w:t></w:r>'
Set-ParagraphXml $p $pPr17_1 $runs

# Insert a brand-new paragraph "Spalte L: ..." right before "Spalte M: ..."
$mIdx = Get-ParaIndexByText $d "*Spalte M: Welcher Gesamtumsatz mit jedem Produkt gemacht wurde.*"
$runs = '<w:r><w:t>Spalte L: Wieviele Einheiten vom Produkt insgesamt bestellt wurden</w:t></w:r>'
Insert-ParagraphBefore $d $mIdx $pPr17_1 $runs

# ---------------------------------------------------------------------
# Change 7: "Nutzen Sie dazu die Funktionen ZÄHLENWENN() und SUMMEWENN()"
# -> add proofErr markers around ZÄHLENWENN(
# ---------------------------------------------------------------------
$p = Get-ParaByText $d "*Nutzen Sie dazu die Funktionen*"
$runs = '<w:r><w:t xml:space="preserve">Nutzen Sie dazu die Funktionen </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>ZÄHLENWENN(</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>) und SUMMEWENN()</w:t></w:r>'
Set-ParagraphXml $p $pPr17_0 $runs

# ---------------------------------------------------------------------
# Change 8: "P4: Nutzen Sie die Funktion SUMMENPRODUKT() auf geeignete Weise."
# -> add proofErr markers around SUMMENPRODUKT(
# ---------------------------------------------------------------------
$p = Get-ParaByText $d "*P4: Nutzen Sie die Funktion*"
$runs = '<w:r><w:t xml:space="preserve">P4: Nutzen Sie die Funktion </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>SUMMENPRODUKT(</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>) auf geeignete Weise.</w:t></w:r>'
Set-ParagraphXml $p $pPr17_1 $runs

# ---------------------------------------------------------------------
# Change 9: "In den drei Tabellen ..." -> move lastRenderedPageBreak to
# the first run and add proofErr markers around SUMMEWENNS(
# ---------------------------------------------------------------------
$pPr18_0 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr>'
$p = Get-ParaByText $d "*In den drei Tabellen*Verkaufsgebiet und Mitarbeiter*"
$runs = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">In den drei Tabellen &quot;Verkäufe pro Produkt und Verkaufsgebiet&quot;, &quot;Umsatz pro Produkt und Verkaufsgebiet&quot; sowie &quot;Umsatz pro Produkt und Verkaufsgebiet und Mitarbeiter&quot; nutzen Sie die Formeln </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>SUMMEWENNS(</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>) bzw. ZÄHLENWENNS() um die gewünschten Werte zu ermitteln. Passen Sie die Zelle V2 an, und kontrollieren Sie, ob sich die Werte in der dritten Tabelle ändern.</w:t></w:r>'
Set-ParagraphXml $p $pPr18_0 $runs
